$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# New shared string value used for rows 8-21 (note lowercase "data")
$newName = "callNumber10AssetNumberS20M20L10.Rdata"

# Column A: widen column, update cells A8:A21 to the new filename text (style untouched by value set)
$ws.Columns.Item(1).ColumnWidth = 43.85546875

for ($r = 8; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = $newName
}

# Apply style index 2 (the "XLConnect.String" right-aligned wrap style used elsewhere in col A)
# to A2:A21 (rows 12-16 already use it; this brings A2:A11 and A17:A21 in line)
$ws.Range("A2:A21").Style = $ws.Range("A12").Style

# H2 becomes a plain number (10) instead of a text value "10"
$ws.Cells.Item(2, 8).Value = 10

# Update the active selection to A23
$ws.Range("A23").Select()

